$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared by every data row (columns K through P)
$K = 2
$L = 0.6666666666666666
$M = 0.2786473333333334
$N = 0.8359420000000001
$O = 1
$P = 1

# Row 2: ECs -> Ccl3/Ackr2 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl3"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.515328
$ws.Range("H2").Value = 1.030656
$ws.Range("I2").Value = 0.0001613041071716636
$ws.Range("J2").Value = 0.0001075508877125494
$ws.Range("K2").Value = $K
$ws.Range("L2").Value = $L
$ws.Range("M2").Value = $M
$ws.Range("N2").Value = $N
$ws.Range("O2").Value = $O
$ws.Range("P2").Value = $P
$ws.Range("Q2").Value = 0.143594772992
$ws.Range("R2").Value = 0.8615686379520001
$ws.Range("S2").Value = 0.0001613041071716636
$ws.Range("T2").Value = 0.0001075508877125494

# Row 3: Inflammatory-Mac -> Ccl3/Ackr2 -> FAPs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Ccl3"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 178.0773926666667
$ws.Range("H3").Value = 534.232178
$ws.Range("I3").Value = 0.05574045041518104
$ws.Range("J3").Value = 0.0557481303058525
$ws.Range("K3").Value = $K
$ws.Range("L3").Value = $L
$ws.Range("M3").Value = $M
$ws.Range("N3").Value = $N
$ws.Range("O3").Value = $O
$ws.Range("P3").Value = $P
$ws.Range("Q3").Value = 49.62079059351956
$ws.Range("R3").Value = 446.587115341676
$ws.Range("S3").Value = 0.05574045041518104
$ws.Range("T3").Value = 0.0557481303058525

# Row 4: MuSCs -> Ccl3/Ackr2 -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ccl3"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.8050075
$ws.Range("H4").Value = 1.610015
$ws.Range("I4").Value = 0.0002519774125488872
$ws.Range("J4").Value = 0.0001680080865783736
$ws.Range("K4").Value = $K
$ws.Range("L4").Value = $L
$ws.Range("M4").Value = $M
$ws.Range("N4").Value = $N
$ws.Range("O4").Value = $O
$ws.Range("P4").Value = $P
$ws.Range("Q4").Value = 0.2243131931883333
$ws.Range("R4").Value = 1.34587915913
$ws.Range("S4").Value = 0.0002519774125488872
$ws.Range("T4").Value = 0.0001680080865783736

# Row 5: Neutrophils -> Ccl3/Ackr2 -> FAPs
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Ccl3"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2891.564290333333
$ws.Range("H5").Value = 8674.692870999999
$ws.Range("I5").Value = 0.9050957762467463
$ws.Range("J5").Value = 0.905220479878615
$ws.Range("K5").Value = $K
$ws.Range("L5").Value = $L
$ws.Range("M5").Value = $M
$ws.Range("N5").Value = $N
$ws.Range("O5").Value = $O
$ws.Range("P5").Value = $P
$ws.Range("Q5").Value = 805.7266786632758
$ws.Range("R5").Value = 7251.540107969482
$ws.Range("S5").Value = 0.9050957762467463
$ws.Range("T5").Value = 0.905220479878615

# Row 6 (new): Resolving-Mac -> Ccl3/Ackr2 -> FAPs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Ccl3"
$ws.Range("C6").Value = "Ackr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 123.798543
$ws.Range("H6").Value = 371.395629
$ws.Range("I6").Value = 0.03875049181835219
$ws.Range("J6").Value = 0.03875583084124156
$ws.Range("K6").Value = $K
$ws.Range("L6").Value = $L
$ws.Range("M6").Value = $M
$ws.Range("N6").Value = $N
$ws.Range("O6").Value = $O
$ws.Range("P6").Value = $P
$ws.Range("Q6").Value = 34.496133877502
$ws.Range("R6").Value = 310.465204897518
$ws.Range("S6").Value = 0.03875049181835219
$ws.Range("T6").Value = 0.03875583084124156
